$wb = $excel.ActiveWorkbook

# --- Sheet "FEINmismatch": update Date column (B) for rows 2-13 ---
$ws1 = $wb.Worksheets.Item("FEINmismatch")
$dates1 = @(
    "Wed Nov 01 15:32:04 EDT 2023",
    "Wed Nov 01 15:32:17 EDT 2023",
    "Wed Nov 01 15:32:28 EDT 2023",
    "Wed Nov 01 15:32:40 EDT 2023",
    "Wed Nov 01 15:32:51 EDT 2023",
    "Wed Nov 01 15:33:03 EDT 2023",
    "Wed Nov 01 15:33:14 EDT 2023",
    "Wed Nov 01 15:33:25 EDT 2023",
    "Wed Nov 01 15:33:36 EDT 2023",
    "Wed Nov 01 15:33:46 EDT 2023",
    "Wed Nov 01 15:33:57 EDT 2023",
    "Wed Nov 01 15:34:08 EDT 2023"
)
for ($i = 0; $i -lt $dates1.Length; $i++) {
    $row = $i + 2
    $ws1.Range("B$row").Value = $dates1[$i]
}

# --- Sheet "FEINSSNmismatch": update Date column (B) for rows 2-17 ---
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$dates2 = @(
    "Wed Nov 01 15:34:19 EDT 2023",
    "Wed Nov 01 15:34:29 EDT 2023",
    "Wed Nov 01 15:34:40 EDT 2023",
    "Wed Nov 01 15:34:50 EDT 2023",
    "Wed Nov 01 15:35:00 EDT 2023",
    "Wed Nov 01 15:35:11 EDT 2023",
    "Wed Nov 01 15:35:21 EDT 2023",
    "Wed Nov 01 15:35:32 EDT 2023",
    "Wed Nov 01 15:35:42 EDT 2023",
    "Wed Nov 01 15:35:53 EDT 2023",
    "Wed Nov 01 15:36:03 EDT 2023",
    "Wed Nov 01 15:36:14 EDT 2023",
    "Wed Nov 01 15:36:24 EDT 2023",
    "Wed Nov 01 15:36:35 EDT 2023",
    "Wed Nov 01 15:36:45 EDT 2023",
    "Wed Nov 01 15:36:55 EDT 2023"
)
for ($i = 0; $i -lt $dates2.Length; $i++) {
    $row = $i + 2
    $ws2.Range("B$row").Value = $dates2[$i]
}
